$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.728.70"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.17"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.031"
$ws.Range("E4").Value = "  +3.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.67"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.029"
$ws.Range("E6").Value = "  +2.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5189"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4001"
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08427"
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.125"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.47"
$ws.Range("E11").Value = "  +2.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.317"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.72"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.302"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.030"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.805.20"
$ws.Range("E16").Value = "  -3.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001122"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.04"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06835"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.91"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.029"
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.036"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.749.47"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.312"
$ws.Range("E25").Value = "  +2.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.47"
$ws.Range("E26").Value = "  +2.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.94"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.014.48"
$ws.Range("E28").Value = "  -3.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.407"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.76"
$ws.Range("E30").Value = "  +2.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1063"
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.893"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.673"
$ws.Range("E34").Value = "  +2.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02458"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06573"
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2206"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.307"
$ws.Range("E38").Value = "  +7.62%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.021"
$ws.Range("E39").Value = "  -5.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.202"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6509"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.064"
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.36"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6099"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.25"
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.771"
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.236"
$ws.Range("E47").Value = "  -3.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.023"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.223"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.76"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06898"
$ws.Range("E51").Value = "  +0.28%  "
